$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the dataset. It belongs right above the
# current row 42, so push everything from row 42 down by one (this also
# grows the used range from R66 to R67) and then fill the freshly opened
# row 42 with the new record's values.
$ws.Rows("42:42").Insert()

$ws.Range("A42").Value = 5
$ws.Range("B42").Value = "Macroferia Regional de Talca"
$ws.Range("C42").Value = "Maule"
$ws.Range("D42").Value = 44488
$ws.Range("E42").Value = 7
$ws.Range("F42").Value = 100112013
$ws.Range("G42").Value = "Alcachofa"
$ws.Range("H42").Value = "Madrigal"
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 300
$ws.Range("K42").Value = 10000
$ws.Range("L42").Value = 10000
$ws.Range("M42").Value = 10000
$ws.Range("N42").Value = "`$/caja 40 unidades"
$ws.Range("O42").Value = "Provincia del Elquí"
$ws.Range("P42").Value = 250
$ws.Range("Q42").Value = 40
$ws.Range("R42").Value = "Hortaliza"
